$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.65
$ws.Range("A7").Value = -21.153
$ws.Range("B7").Value = 6.053
$ws.Range("B12").Value = 5.401999999999999
$ws.Range("D13").Value = -8.281000000000001
$ws.Range("D14").Value = -7.786999999999999
$ws.Range("B15").Value = 5.359000000000001
$ws.Range("A16").Value = -21.43
$ws.Range("D16").Value = -8.518000000000001
$ws.Range("D19").Value = -7.796000000000001
$ws.Range("A20").Value = -21.857
$ws.Range("B20").Value = 5.539
$ws.Range("B21").Value = 8.904
$ws.Range("B22").Value = 6.4
$ws.Range("D22").Value = -8.377000000000001
$ws.Range("B23").Value = 7.359999999999999
$ws.Range("A28").Value = -21.878
$ws.Range("A29").Value = -21.675
$ws.Range("B29").Value = 5.731
$ws.Range("A32").Value = -21.653
$ws.Range("B34").Value = 8.059000000000001
$ws.Range("D36").Value = -7.834000000000001
$ws.Range("A40").Value = -20.312
$ws.Range("B42").Value = 7.238
$ws.Range("B43").Value = 5.529000000000001
$ws.Range("B44").Value = 5.231
$ws.Range("B45").Value = 5.286
$ws.Range("A46").Value = -20.849
$ws.Range("B46").Value = 6.865
$ws.Range("D46").Value = -8.044999999999998
$ws.Range("B50").Value = 5.828
$ws.Range("D50").Value = -8.158000000000001
$ws.Range("A51").Value = -20.771
$ws.Range("B51").Value = 7.779000000000001
$ws.Range("A52").Value = -21.316
$ws.Range("A57").Value = -22.137
$ws.Range("A59").Value = -22.091
$ws.Range("A62").Value = -21.85
$ws.Range("A66").Value = -21.44600000000001
$ws.Range("B66").Value = 5.709000000000001
$ws.Range("B67").Value = 5.194999999999999
$ws.Range("A73").Value = -20.53
$ws.Range("A74").Value = -21.043
$ws.Range("B79").Value = 5.680999999999999
$ws.Range("B84").Value = 5.781000000000001
$ws.Range("A92").Value = -21.333
$ws.Range("B92").Value = 5.527
$ws.Range("D95").Value = -7.815
$ws.Range("B97").Value = 5.161
$ws.Range("D97").Value = -8.465
$ws.Range("A100").Value = -21.481
